# Apply the Class04.pptx edit:
#  1. On slide 1, remove the "Research Ethics Presentations" bullet and merge the
#     "Writing " + "Genres" runs into a single "Writing Genres" run.
#  2. Delete the two duplicate "Research Ethics (Team presentations)" slides
#     (originally slides 2 and 3).

$p = $ppt.ActivePresentation

# --- Step 1: edit slide 1's bullet list textbox ---
$slide1 = $p.Slides.Item(1)
$bulletShape = $slide1.Shapes.Item(2)
$tr = $bulletShape.TextFrame.TextRange

# Remove the "Research Ethics Presentations" paragraph entirely (it's paragraph 3:
# "Class #4", "Faculty Presentations", "Research Ethics Presentations", "Writing "+"Genres", "Writing in Physics").
$tr.Paragraphs(3, 1).Delete()

# The former paragraph 4 ("Writing " + "Genres") is now paragraph 3.
$genresPara = $tr.Paragraphs(3, 1)

# Remove the leading "Writing " run text (first 8 characters, including trailing space)
# so only the second run ("Genres", dirty="0") remains.
$genresPara.Characters(1, 8).Text = ""

# Prepend "Writing " back onto the remaining run so the whole phrase reads "Writing Genres"
# while keeping that run's original formatting (single merged run).
$genresPara.Runs(1).Text = "Writing Genres"

# --- Step 2: remove the two duplicate "Research Ethics" slides (now slides 2 and 3) ---
$p.Slides.Item(3).Delete()
$p.Slides.Item(2).Delete()
